$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Professional summary paragraph: simple text replace (no formatting change)
#    "... errors affecting all Black and Asian-American voters, developed ..."
#    -> "... errors affecting 50M voters, developed ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "errors affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "errors affecting 50M voters, developed geospatial ML", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Work-experience bullet under "Partner - Siege Analytics": the phrase
#    "all Black and Asian-American" becomes a bold/colored "50M" run, while
#    the surrounding text (incl. " voters, ...") stays in plain runs, i.e.
#    the single run gets split into three runs.
# ---------------------------------------------------------------------------
$bulletParaIdx = Get-ParaIndexByText $d "Discovered systematic race coding errors affecting"
$bulletPara = $d.Paragraphs.Item($bulletParaIdx)
$bulletRange = $bulletPara.Range.Duplicate
$bulletRange.Find.Execute("all Black and Asian-American") | Out-Null
$bulletRange.Text = "50M"
$bulletRange.Font.Bold = 1
$rr = 0x2C
$gg = 0x3E
$bb = 0x50
$bulletRange.Font.Color = ($bb * 65536) + ($gg * 256) + $rr

# ---------------------------------------------------------------------------
# 3. Reorder work-experience blocks: move the "Analytics Supervisor - GSD&M"
#    entry (heading + 4 body paragraphs = 5 paragraphs total) so that it
#    follows the "Data Products Manager - Helm/Murmuration" entry (also 5
#    paragraphs) instead of preceding it.
# ---------------------------------------------------------------------------
$gsdmHeadingIdx = Get-ParaIndexByText $d "Analytics Supervisor - GSD&M"
$gsdmStart = $d.Paragraphs.Item($gsdmHeadingIdx)
$gsdmEnd = $d.Paragraphs.Item($gsdmHeadingIdx + 4)
$gsdmBlock = $d.Range($gsdmStart.Range.Start, $gsdmEnd.Range.End)
$gsdmBlock.Cut()

$helmHeadingIdx = Get-ParaIndexByText $d "Data Products Manager - Helm/Murmuration"
$helmEnd = $d.Paragraphs.Item($helmHeadingIdx + 4)
$pasteRange = $d.Range($helmEnd.Range.End, $helmEnd.Range.End)
$pasteRange.Paste()

# Cut/Paste that lands exactly on a paragraph-mark boundary can drop the
# heading style of the first pasted paragraph; restore it explicitly.
$restoredHeadingIdx = Get-ParaIndexByText $d "Analytics Supervisor - GSD&M"
$d.Paragraphs.Item($restoredHeadingIdx).Style = "Heading 3"

# ---------------------------------------------------------------------------
# 4. "Geospatial Demographic Classification System" project impact line:
#    simple text replace that also appends "nationwide".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "data affecting all Black and Asian-American voters, improved electoral",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "data affecting 50M voters nationwide, improved electoral", 2) | Out-Null
